$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24
$ws.Cells.Item($row, 1).Value = "2025/12/03 16:00"
$ws.Cells.Item($row, 2).Value = "-"
$ws.Cells.Item($row, 3).Value = "-"
$ws.Cells.Item($row, 4).Value = "-"
$ws.Cells.Item($row, 5).Value = "-"
$ws.Cells.Item($row, 6).Value = "-"
$ws.Cells.Item($row, 7).Value = "-"
